$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 34.111822
$ws.Range("H2").Value = 102.335466
$ws.Range("I2").Value = 0.4228853893909983
$ws.Range("J2").Value = 0.4228853893909983
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 162.98837
$ws.Range("N2").Value = 488.96511
$ws.Range("O2").Value = 0.9909539753179891
$ws.Range("P2").Value = 0.9909539753179891
$ws.Range("Q2").Value = 5559.83026551014
$ws.Range("R2").Value = 50038.47238959125
$ws.Range("S2").Value = 0.4190599577209056
$ws.Range("T2").Value = 0.4190599577209056
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 34.111822
$ws.Range("H3").Value = 102.335466
$ws.Range("I3").Value = 0.4228853893909983
$ws.Range("J3").Value = 0.4228853893909983
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.6513563333333333
$ws.Range("N3").Value = 1.954069
$ws.Range("O3").Value = 0.003960185305646138
$ws.Range("P3").Value = 0.003960185305646138
$ws.Range("Q3").Value = 22.21895130123933
$ws.Range("R3").Value = 199.970561711154
$ws.Range("S3").Value = 0.001674704505038677
$ws.Range("T3").Value = 0.001674704505038677
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 34.111822
$ws.Range("H4").Value = 102.335466
$ws.Range("I4").Value = 0.4228853893909983
$ws.Range("J4").Value = 0.4228853893909983
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.8364996666666666
$ws.Range("N4").Value = 2.509499
$ws.Range("O4").Value = 0.005085839376364744
$ws.Range("P4").Value = 0.005085839376364744
$ws.Range("Q4").Value = 28.53452773239266
$ws.Range("R4").Value = 256.810749591534
$ws.Range("S4").Value = 0.002150727165054077
$ws.Range("T4").Value = 0.002150727165054077
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 34.88211266666666
$ws.Range("H5").Value = 104.646338
$ws.Range("I5").Value = 0.4324347083490296
$ws.Range("J5").Value = 0.4324347083490295
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 162.98837
$ws.Range("N5").Value = 488.96511
$ws.Range("O5").Value = 0.9909539753179891
$ws.Range("P5").Value = 0.9909539753179891
$ws.Range("Q5").Value = 5685.378685696353
$ws.Range("R5").Value = 51168.40817126718
$ws.Range("S5").Value = 0.4285228933039461
$ws.Range("T5").Value = 0.4285228933039461
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 34.88211266666666
$ws.Range("H6").Value = 104.646338
$ws.Range("I6").Value = 0.4324347083490296
$ws.Range("J6").Value = 0.4324347083490295
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 0.6513563333333333
$ws.Range("N6").Value = 1.954069
$ws.Range("O6").Value = 0.003960185305646138
$ws.Range("P6").Value = 0.003960185305646138
$ws.Range("Q6").Value = 22.72068500548022
$ws.Range("R6").Value = 204.486165049322
$ws.Range("S6").Value = 0.0017125215776552
$ws.Range("T6").Value = 0.0017125215776552
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 34.88211266666666
$ws.Range("H7").Value = 104.646338
$ws.Range("I7").Value = 0.4324347083490296
$ws.Range("J7").Value = 0.4324347083490295
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.8364996666666666
$ws.Range("N7").Value = 2.509499
$ws.Range("O7").Value = 0.005085839376364744
$ws.Range("P7").Value = 0.005085839376364744
$ws.Range("Q7").Value = 29.17887561829577
$ws.Range("R7").Value = 262.609880564662
$ws.Range("S7").Value = 0.002199293467428298
$ws.Range("T7").Value = 0.002199293467428298
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 11.67052633333333
$ws.Range("H8").Value = 35.011579
$ws.Range("I8").Value = 0.1446799022599722
$ws.Range("J8").Value = 0.1446799022599721
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 162.98837
$ws.Range("N8").Value = 488.96511
$ws.Range("O8").Value = 0.9909539753179891
$ws.Range("P8").Value = 0.9909539753179891
$ws.Range("Q8").Value = 1902.160064112077
$ws.Range("R8").Value = 17119.44057700869
$ws.Range("S8").Value = 0.1433711242931375
$ws.Range("T8").Value = 0.1433711242931375
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 11.67052633333333
$ws.Range("H9").Value = 35.011579
$ws.Range("I9").Value = 0.1446799022599722
$ws.Range("J9").Value = 0.1446799022599721
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.6513563333333333
$ws.Range("N9").Value = 1.954069
$ws.Range("O9").Value = 0.003960185305646138
$ws.Range("P9").Value = 0.003960185305646138
$ws.Range("Q9").Value = 7.601671240550111
$ws.Range("R9").Value = 68.41504116495099
$ws.Range("S9").Value = 0.0005729592229522612
$ws.Range("T9").Value = 0.0005729592229522611
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 11.67052633333333
$ws.Range("H10").Value = 35.011579
$ws.Range("I10").Value = 0.1446799022599722
$ws.Range("J10").Value = 0.1446799022599721
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.8364996666666666
$ws.Range("N10").Value = 2.509499
$ws.Range("O10").Value = 0.005085839376364744
$ws.Range("P10").Value = 0.005085839376364744
$ws.Range("Q10").Value = 9.762391387657889
$ws.Range("R10").Value = 87.861522488921
$ws.Range("S10").Value = 0.0007358187438823689
$ws.Range("T10").Value = 0.0007358187438823687
